$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder "42*30.44" formulas in column F with the
# patient-specific timepoint day counts (hard values, no formula).
$timepointDays = @{
    6  = 386
    17 = 365
    20 = 2449
    23 = 903
    26 = 2534
    29 = 2197
    32 = 1567
    35 = 2205
    38 = 1048
    41 = 1138
    44 = 945
    47 = 754
    50 = 2345
    53 = 1099
    59 = 880
    62 = 951
}

foreach ($row in $timepointDays.Keys) {
    $ws.Range("F$row").Value = $timepointDays[$row]
}

# Update the saved view/selection: scroll back to the top and select F11.
$ws.Range("F11").Select()
